$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 old worker detail rows (17-21), keeping only the first
# worker row (16); this also drops the now-unused shared strings for
# those workers and renumbers the trailing signature rows 26/27 -> 21/22.
$ws.Rows("17:21").Delete()

# Update the "VALOR MORA" total for the remaining worker.
$ws.Range("E11").Value = 40000

# Update worker / period counters.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D's best-fit width shrinks now that the longest name
# ("SCARLETT STHEFANI CASTRO TAJAN") was removed along with its row.
$ws.Range("D1").ColumnWidth = 24.83
